$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.627.00"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.895.33"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.67"
$ws.Range("E5").Value = "  +2.67%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  +1.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2950"
$ws.Range("E8").Value = "  +2.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06728"
$ws.Range("E9").Value = "  +1.51%  "
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.18"
$ws.Range("E10").Value = "  +3.19%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.867.84"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07359"
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.151"
$ws.Range("E13").Value = "  +3.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.40"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6709"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").Value = "30.567.20"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007885"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.47"
$ws.Range("E18").Value = "  +4.65%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "2.132.45"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.305"
$ws.Range("E21").Value = "  +12.51%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "191.65"
$ws.Range("E23").Value = "  +2.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.226"
$ws.Range("E24").Value = "  +3.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.562"
$ws.Range("E25").Value = "  +3.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.89"
$ws.Range("E26").Value = "  +3.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.46"
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.946"
$ws.Range("E28").Value = "  +6.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.475"
$ws.Range("E29").Value = "  +5.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.453"
$ws.Range("E30").Value = "  +5.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09199"
$ws.Range("E31").Value = "  +2.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.193"
$ws.Range("E32").Value = "  +7.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05244"
$ws.Range("E33").Value = "  +1.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7478"
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.107"
$ws.Range("E35").Value = "  +3.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.714"
$ws.Range("E36").Value = "  +0.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01837"
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.695"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9232"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.061"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4429"
$ws.Range("E41").Value = "  +3.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.962"
$ws.Range("E42").Value = "  +4.89%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.30"
$ws.Range("E43").Value = "  +2.99%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.23"
$ws.Range("E44").Value = "  +25.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9941"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1388"
$ws.Range("E46").Value = "  +4.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.596"
$ws.Range("E47").Value = "  +5.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.080"
$ws.Range("E48").Value = "  +5.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.10"
$ws.Range("E49").Value = "  +6.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05831"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3962"
$ws.Range("E51").Value = "  +2.46%  "
